$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained 23 new incident rows (rows 73-95) and every existing RISP label
# was renamed from "<n>º RISP - <city>" to "<n>º Departamento - <city>".
# Extend formatting (styles/number formats/borders) from the last existing data row
# down across the newly added rows before writing values into them.
$ws.Range("A72:I72").Copy()
$ws.Range("A73:I95").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Full target data for rows 2-95: MUNICIPIO_COD, MUNICIPIO_FATO, DATA_FATO (serial),
# MES, ANO, RISP, RMBH, TENTADO/CONSUMADO, QTDE_VITIMAS
$data = @(
  @(310620, "BELO HORIZONTE", 44927, 1, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(311570, "CENTRAL DE MINAS", 44927, 1, 2023, "8º Departamento - Governador Valadares", "3) Interior de MG", "TENTADO", 1),
  @(312245, "DIVISOPOLIS", 44927, 1, 2023, "15º Departamento - Teófilo Otoni", "3) Interior de MG", "CONSUMADO", 1),
  @(310350, "ARAGUARI", 44928, 1, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "CONSUMADO", 1),
  @(315180, "POCOS DE CALDAS", 44928, 1, 2023, "18º Departamento - Poços de Caldas", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 44932, 1, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(311860, "CONTAGEM", 44933, 1, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(316870, "TIMOTEO", 44934, 1, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "TENTADO", 1),
  @(311930, "COROMANDEL", 44934, 1, 2023, "10º Departamento - Patos de Minas", "3) Interior de MG", "CONSUMADO", 1),
  @(312260, "DOM JOAQUIM", 44935, 1, 2023, "14º Departamento - Curvelo", "3) Interior de MG", "TENTADO", 1),
  @(314280, "MONTE ALEGRE DE MINAS", 44935, 1, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "CONSUMADO", 1),
  @(314330, "MONTES CLAROS", 44937, 1, 2023, "11º Departamento - Montes Claros", "3) Interior de MG", "TENTADO", 1),
  @(311340, "CARATINGA", 44938, 1, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "TENTADO", 1),
  @(314320, "MONTE SANTO DE MINAS", 44938, 1, 2023, "18º Departamento - Poços de Caldas", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 44940, 1, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(314655, "PAI PEDRO", 44940, 1, 2023, "11º Departamento - Montes Claros", "3) Interior de MG", "TENTADO", 1),
  @(310690, "BICAS", 44941, 1, 2023, "4º Departamento - Juiz de Fora", "3) Interior de MG", "CONSUMADO", 1),
  @(313010, "IGARAPE", 44941, 1, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(313550, "JEQUERI", 44941, 1, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "CONSUMADO", 1),
  @(311570, "CENTRAL DE MINAS", 44941, 1, 2023, "8º Departamento - Governador Valadares", "3) Interior de MG", "CONSUMADO", 1),
  @(317020, "UBERLANDIA", 44947, 1, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "TENTADO", 1),
  @(312430, "ESPINOSA", 44947, 1, 2023, "11º Departamento - Montes Claros", "3) Interior de MG", "CONSUMADO", 1),
  @(313440, "ITURAMA", 44947, 1, 2023, "5º Departamento - Uberaba", "3) Interior de MG", "CONSUMADO", 1),
  @(311070, "CAMBUQUIRA", 44951, 1, 2023, "6º Departamento - Lavras", "3) Interior de MG", "TENTADO", 1),
  @(314140, "MEDINA", 44953, 1, 2023, "15º Departamento - Teófilo Otoni", "3) Interior de MG", "TENTADO", 1),
  @(317040, "UNAI", 44955, 1, 2023, "16º Departamento - Unaí", "3) Interior de MG", "TENTADO", 1),
  @(317040, "UNAI", 44955, 1, 2023, "16º Departamento - Unaí", "3) Interior de MG", "CONSUMADO", 1),
  @(312230, "DIVINOPOLIS", 44956, 1, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 1),
  @(310160, "ALFENAS", 44959, 2, 2023, "18º Departamento - Poços de Caldas", "3) Interior de MG", "TENTADO", 1),
  @(311860, "CONTAGEM", 44959, 2, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(314110, "MATOZINHOS", 44959, 2, 2023, "19º Departamento - Sete Lagoas", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(315950, "SANTA RITA DO ITUETO", 44963, 2, 2023, "8º Departamento - Governador Valadares", "3) Interior de MG", "CONSUMADO", 1),
  @(314930, "PEDRO LEOPOLDO", 44964, 2, 2023, "19º Departamento - Sete Lagoas", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(316940, "TRES PONTAS", 44965, 2, 2023, "6º Departamento - Lavras", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 44968, 2, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "CONSUMADO", 1),
  @(310670, "BETIM", 44969, 2, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(314610, "OURO PRETO", 44971, 2, 2023, "3º Departamento - Vespasiano", "3) Interior de MG", "CONSUMADO", 1),
  @(317020, "UBERLANDIA", 44973, 2, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 44976, 2, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(314465, "NINHEIRA", 44977, 2, 2023, "11º Departamento - Montes Claros", "3) Interior de MG", "CONSUMADO", 1),
  @(312760, "GOUVEIA", 44979, 2, 2023, "14º Departamento - Curvelo", "3) Interior de MG", "TENTADO", 1),
  @(310620, "BELO HORIZONTE", 44980, 2, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 2),
  @(311860, "CONTAGEM", 44982, 2, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(310670, "BETIM", 44982, 2, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(313690, "JURUAIA", 44982, 2, 2023, "18º Departamento - Poços de Caldas", "3) Interior de MG", "CONSUMADO", 1),
  @(313130, "IPATINGA", 44983, 2, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "TENTADO", 1),
  @(313670, "JUIZ DE FORA", 44983, 2, 2023, "4º Departamento - Juiz de Fora", "3) Interior de MG", "CONSUMADO", 1),
  @(315725, "SANTA BARBARA DO LESTE", 44984, 2, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "CONSUMADO", 1),
  @(314390, "MURIAE", 44985, 2, 2023, "4º Departamento - Juiz de Fora", "3) Interior de MG", "TENTADO", 1),
  @(316720, "SETE LAGOAS", 44986, 3, 2023, "19º Departamento - Sete Lagoas", "3) Interior de MG", "TENTADO", 1),
  @(310620, "BELO HORIZONTE", 44987, 3, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "CONSUMADO", 1),
  @(315460, "RIBEIRAO DAS NEVES", 44988, 3, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(317130, "VICOSA", 44988, 3, 2023, "4º Departamento - Juiz de Fora", "3) Interior de MG", "TENTADO", 1),
  @(310670, "BETIM", 44990, 3, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(312160, "DIAMANTINA", 44991, 3, 2023, "14º Departamento - Curvelo", "3) Interior de MG", "CONSUMADO", 1),
  @(310240, "ALVORADA DE MINAS", 44993, 3, 2023, "14º Departamento - Curvelo", "3) Interior de MG", "TENTADO", 1),
  @(313670, "JUIZ DE FORA", 44993, 3, 2023, "4º Departamento - Juiz de Fora", "3) Interior de MG", "TENTADO", 1),
  @(310620, "BELO HORIZONTE", 44996, 3, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "CONSUMADO", 1),
  @(311820, "CONQUISTA", 44997, 3, 2023, "5º Departamento - Uberaba", "3) Interior de MG", "TENTADO", 1),
  @(310620, "BELO HORIZONTE", 45000, 3, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "CONSUMADO", 2),
  @(310720, "BOCAINA DE MINAS", 45000, 3, 2023, "17º Departamento - Pouso Alegre", "3) Interior de MG", "CONSUMADO", 1),
  @(313130, "IPATINGA", 45001, 3, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "CONSUMADO", 1),
  @(314530, "NOVO CRUZEIRO", 45001, 3, 2023, "15º Departamento - Teófilo Otoni", "3) Interior de MG", "CONSUMADO", 1),
  @(315140, "PITANGUI", 45003, 3, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 1),
  @(314630, "PADRE PARAISO", 45003, 3, 2023, "15º Departamento - Teófilo Otoni", "3) Interior de MG", "CONSUMADO", 1),
  @(310180, "ALPERCATA", 45005, 3, 2023, "8º Departamento - Governador Valadares", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 45006, 3, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(315140, "PITANGUI", 45007, 3, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 1),
  @(313190, "ITABIRITO", 45011, 3, 2023, "3º Departamento - Vespasiano", "3) Interior de MG", "TENTADO", 1),
  @(313420, "ITUIUTABA", 45011, 3, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "CONSUMADO", 1),
  @(311940, "CORONEL FABRICIANO", 45011, 3, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "CONSUMADO", 1),
  @(314800, "PATOS DE MINAS", 45012, 3, 2023, "10º Departamento - Patos de Minas", "3) Interior de MG", "TENTADO", 1),
  @(313240, "ITAJUBA", 45017, 4, 2023, "17º Departamento - Pouso Alegre", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 45018, 4, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(315140, "PITANGUI", 45018, 4, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 1),
  @(313190, "ITABIRITO", 45021, 4, 2023, "3º Departamento - Vespasiano", "3) Interior de MG", "TENTADO", 1),
  @(314600, "OURO FINO", 45021, 4, 2023, "17º Departamento - Pouso Alegre", "3) Interior de MG", "CONSUMADO", 1),
  @(310350, "ARAGUARI", 45022, 4, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "TENTADO", 1),
  @(313940, "MANHUACU", 45022, 4, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "TENTADO", 1),
  @(315895, "SANTANA DO PARAISO", 45023, 4, 2023, "12º Departamento - Ipatinga", "3) Interior de MG", "CONSUMADO", 1),
  @(312980, "IBIRITE", 45024, 4, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(311230, "CAPELINHA", 45026, 4, 2023, "14º Departamento - Curvelo", "3) Interior de MG", "TENTADO", 1),
  @(311860, "CONTAGEM", 45030, 4, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "CONSUMADO", 1),
  @(313020, "IGARATINGA", 45031, 4, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 2),
  @(311830, "CONSELHEIRO LAFAIETE", 45031, 4, 2023, "13º Departamento - Barbacena", "3) Interior de MG", "CONSUMADO", 1),
  @(314330, "MONTES CLAROS", 45032, 4, 2023, "11º Departamento - Montes Claros", "3) Interior de MG", "TENTADO", 1),
  @(314710, "PARA DE MINAS", 45033, 4, 2023, "7º Departamento - Divinópolis", "3) Interior de MG", "TENTADO", 1),
  @(314700, "PARACATU", 45038, 4, 2023, "16º Departamento - Unaí", "3) Interior de MG", "CONSUMADO", 1),
  @(312980, "IBIRITE", 45042, 4, 2023, "2º Departamento - Contagem", "2) RMBH (sem BH)", "TENTADO", 1),
  @(315100, "PIRANGUINHO", 45042, 4, 2023, "17º Departamento - Pouso Alegre", "3) Interior de MG", "TENTADO", 1),
  @(311260, "CAPINOPOLIS", 45042, 4, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "CONSUMADO", 1),
  @(310620, "BELO HORIZONTE", 45045, 4, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "TENTADO", 1),
  @(310620, "BELO HORIZONTE", 45045, 4, 2023, "1º Departamento - Belo Horizonte", "1) Belo Horizonte", "CONSUMADO", 1),
  @(317020, "UBERLANDIA", 45046, 4, 2023, "9º Departamento - Uberlândia", "3) Interior de MG", "TENTADO", 1)
)

$startRow = 2
for ($idx = 0; $idx -lt $data.Count; $idx++) {
  $row = $data[$idx]
  $r = $startRow + $idx
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $ws.Cells.Item($r, 9).Value = $row[8]
}

# Match the author's final view state: scrolled near the bottom with the newly
# added rows selected.
$ws.Range("A73:I95").Select()
$excel.ActiveWindow.ScrollRow = 70